{"js": "// Replace the date paragraph and each \"a\u00d7b=c\" answer cell with its\n// updated value. Every old value in `replacements` is unique within the\n// document, so an exact (non-wildcard) search safely targets exactly one\n// run each time.\nconst replacements = [\n  [\"2024-12-09 Monday\", \"2024-12-10 Tuesday\"],\n  [\"21\u00d726=546\", \"54\u00d748=2592\"],\n  [\"65\u00d799=6435\", \"54\u00d745=2430\"],\n  [\"30\u00d726=780\", \"36\u00d796=3456\"],\n  [\"21\u00d793=1953\", \"22\u00d722=484\"],\n  [\"73\u00d780=5840\", \"90\u00d719=1710\"],\n  [\"34\u00d722=748\", \"31\u00d734=1054\"],\n  [\"95\u00d734=3230\", \"48\u00d782=3936\"],\n  [\"48\u00d759=2832\", \"97\u00d711=1067\"],\n  [\"66\u00d742=2772\", \"75\u00d759=4425\"],\n  [\"22\u00d743=946\", \"50\u00d753=2650\"],\n  [\"39\u00d799=3861\", \"41\u00d798=4018\"],\n  [\"92\u00d733=3036\", \"13\u00d780=1040\"],\n  [\"95\u00d747=4465\", \"40\u00d744=1760\"],\n  [\"28\u00d763=1764\", \"31\u00d786=2666\"],\n  [\"64\u00d745=2880\", \"90\u00d739=3510\"],\n  [\"65\u00d739=2535\", \"61\u00d737=2257\"],\n  [\"78\u00d739=3042\", \"33\u00d785=2805\"],\n  [\"90\u00d781=7290\", \"20\u00d718=360\"],\n  [\"87\u00d740=3480\", \"17\u00d767=1139\"],\n  [\"84\u00d713=1092\", \"59\u00d782=4838\"],\n  [\"13\u00d742=546\", \"82\u00d739=3198\"],\n  [\"93\u00d752=4836\", \"20\u00d796=1920\"],\n  [\"20\u00d769=1380\", \"30\u00d743=1290\"],\n  [\"90\u00d755=4950\", \"52\u00d762=3224\"],\n  [\"54\u00d787=4698\", \"92\u00d713=1196\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Replace the date paragraph and each \"a\u00d7b=c\" answer cell with its\n# updated value. Every old value is unique within the document, so a\n# plain (non-wildcard) Find/Replace targets exactly one run each time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-12-09 Monday', '2024-12-10 Tuesday'),\n    @('21\u00d726=546', '54\u00d748=2592'),\n    @('65\u00d799=6435', '54\u00d745=2430'),\n    @('30\u00d726=780', '36\u00d796=3456'),\n    @('21\u00d793=1953', '22\u00d722=484'),\n    @('73\u00d780=5840', '90\u00d719=1710'),\n    @('34\u00d722=748', '31\u00d734=1054'),\n    @('95\u00d734=3230', '48\u00d782=3936'),\n    @('48\u00d759=2832', '97\u00d711=1067'),\n    @('66\u00d742=2772', '75\u00d759=4425'),\n    @('22\u00d743=946', '50\u00d753=2650'),\n    @('39\u00d799=3861', '41\u00d798=4018'),\n    @('92\u00d733=3036', '13\u00d780=1040'),\n    @('95\u00d747=4465', '40\u00d744=1760'),\n    @('28\u00d763=1764', '31\u00d786=2666'),\n    @('64\u00d745=2880', '90\u00d739=3510'),\n    @('65\u00d739=2535', '61\u00d737=2257'),\n    @('78\u00d739=3042', '33\u00d785=2805'),\n    @('90\u00d781=7290', '20\u00d718=360'),\n    @('87\u00d740=3480', '17\u00d767=1139'),\n    @('84\u00d713=1092', '59\u00d782=4838'),\n    @('13\u00d742=546', '82\u00d739=3198'),\n    @('93\u00d752=4836', '20\u00d796=1920'),\n    @('20\u00d769=1380', '30\u00d743=1290'),\n    @('90\u00d755=4950', '52\u00d762=3224'),\n    @('54\u00d787=4698', '92\u00d713=1196'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
